# Auto-refresh BRVM recommendations (GitHub Actions job)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recommandations")

# --- Rows 2-15: update "Jours en Baisse" (C) and "Variation Totale (%)" (D) ---
$ws.Cells.Item(2, 3).Value = 6
$ws.Cells.Item(2, 4).Value = 2519.11
$ws.Cells.Item(3, 3).Value = 3
$ws.Cells.Item(3, 4).Value = 2025
$ws.Cells.Item(4, 3).Value = 3
$ws.Cells.Item(4, 4).Value = 1945
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 4).Value = 1807.54
$ws.Cells.Item(6, 3).Value = 3
$ws.Cells.Item(6, 4).Value = 1416.82
$ws.Cells.Item(7, 3).Value = 3
$ws.Cells.Item(7, 4).Value = 1127.44
$ws.Cells.Item(8, 3).Value = 3
$ws.Cells.Item(8, 4).Value = 1083.12
$ws.Cells.Item(9, 3).Value = 3
$ws.Cells.Item(9, 4).Value = 475.28
$ws.Cells.Item(10, 3).Value = 3
$ws.Cells.Item(10, 4).Value = 411.1
$ws.Cells.Item(11, 3).Value = 3
$ws.Cells.Item(11, 4).Value = 406.43
$ws.Cells.Item(12, 3).Value = 3
$ws.Cells.Item(12, 4).Value = 399.44
$ws.Cells.Item(13, 3).Value = 3
$ws.Cells.Item(13, 4).Value = 388.04
$ws.Cells.Item(14, 3).Value = 3
$ws.Cells.Item(14, 4).Value = 329.76
$ws.Cells.Item(15, 3).Value = 3
$ws.Cells.Item(15, 4).Value = 287.09

# --- Rows 22-42: full refresh of recommendation table (re-ranked) ---
$ws.Cells.Item(22, 1).Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$ws.Cells.Item(22, 2).Value = 3
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = 22.16
$ws.Cells.Item(22, 5).Value = 7.37
$ws.Cells.Item(22, 6).Value = '🟢 Achat'
$ws.Cells.Item(22, 7).Value = '✅ Renforcer'

$ws.Cells.Item(23, 1).Value = 'SICABLE CI (CABC)'
$ws.Cells.Item(23, 2).Value = 2
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 14.78
$ws.Cells.Item(23, 5).Value = 7.45
$ws.Cells.Item(23, 6).Value = '🟡 Observer'
$ws.Cells.Item(23, 7).Value = '➖ Neutre'

$ws.Cells.Item(24, 1).Value = 'SETAO CI (STAC)'
$ws.Cells.Item(24, 2).Value = 2
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 14.59
$ws.Cells.Item(24, 5).Value = 7.27
$ws.Cells.Item(24, 6).Value = '🟡 Observer'
$ws.Cells.Item(24, 7).Value = '➖ Neutre'

$ws.Cells.Item(25, 1).Value = 'CFAO MOTORS CI (CFAC)'
$ws.Cells.Item(25, 2).Value = 1
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 7.35
$ws.Cells.Item(25, 5).Value = 7.35
$ws.Cells.Item(25, 6).Value = '🟡 Observer'
$ws.Cells.Item(25, 7).Value = '➖ Neutre'

$ws.Cells.Item(26, 1).Value = 'SAFCA CI (SAFC)'
$ws.Cells.Item(26, 2).Value = 1
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 7.34
$ws.Cells.Item(26, 5).Value = 7.34
$ws.Cells.Item(26, 6).Value = '🟡 Observer'
$ws.Cells.Item(26, 7).Value = '➖ Neutre'

$ws.Cells.Item(27, 1).Value = 'FILTISAC CI (FTSC)'
$ws.Cells.Item(27, 2).Value = 2
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(27, 4).Value = 7.31
$ws.Cells.Item(27, 5).Value = -7.46
$ws.Cells.Item(27, 6).Value = '🟡 Observer'
$ws.Cells.Item(27, 7).Value = '👀 À surveiller'

$ws.Cells.Item(28, 1).Value = 'SUCRIVOIRE (SCRC)'
$ws.Cells.Item(28, 2).Value = 1
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 6.78
$ws.Cells.Item(28, 5).Value = 6.78
$ws.Cells.Item(28, 6).Value = '🟡 Observer'
$ws.Cells.Item(28, 7).Value = '➖ Neutre'

$ws.Cells.Item(29, 1).Value = 'ONATEL BF (ONTBF)'
$ws.Cells.Item(29, 2).Value = 1
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(29, 4).Value = 5.41
$ws.Cells.Item(29, 5).Value = 7.31
$ws.Cells.Item(29, 6).Value = '🟡 Observer'
$ws.Cells.Item(29, 7).Value = '👀 À surveiller'

$ws.Cells.Item(30, 1).Value = 'BERNABE CI (BNBC)'
$ws.Cells.Item(30, 2).Value = 1
$ws.Cells.Item(30, 3).Value = 1
$ws.Cells.Item(30, 4).Value = 4.34
$ws.Cells.Item(30, 5).Value = -3.16
$ws.Cells.Item(30, 6).Value = '🟡 Observer'
$ws.Cells.Item(30, 7).Value = '👀 À surveiller'

$ws.Cells.Item(31, 1).Value = 'SAPH CI (SPHC)'
$ws.Cells.Item(31, 2).Value = 1
$ws.Cells.Item(31, 3).Value = 1
$ws.Cells.Item(31, 4).Value = 0.64
$ws.Cells.Item(31, 5).Value = -6.81
$ws.Cells.Item(31, 6).Value = '🟡 Observer'
$ws.Cells.Item(31, 7).Value = '👀 À surveiller'

$ws.Cells.Item(32, 1).Value = 'TOTAL'
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(32, 3).Value = 4
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 6).Value = '🟡 Observer'
$ws.Cells.Item(32, 7).Value = '➖ Neutre'

$ws.Cells.Item(33, 1).Value = 'BANK OF AFRICA ML (BOAM)'
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(33, 4).Value = -1.15
$ws.Cells.Item(33, 5).Value = -1.15
$ws.Cells.Item(33, 6).Value = '🟡 Observer'
$ws.Cells.Item(33, 7).Value = '➖ Neutre'

$ws.Cells.Item(34, 1).Value = 'SOGB CI (SOGC)'
$ws.Cells.Item(34, 2).Value = 0
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 4).Value = -1.18
$ws.Cells.Item(34, 5).Value = -1.18
$ws.Cells.Item(34, 6).Value = '🟡 Observer'
$ws.Cells.Item(34, 7).Value = '➖ Neutre'

$ws.Cells.Item(35, 1).Value = 'NSIA BANQUE COTE D''IVOIRE (NSBC)'
$ws.Cells.Item(35, 2).Value = 0
$ws.Cells.Item(35, 3).Value = 1
$ws.Cells.Item(35, 4).Value = -1.83
$ws.Cells.Item(35, 5).Value = -1.83
$ws.Cells.Item(35, 6).Value = '🟡 Observer'
$ws.Cells.Item(35, 7).Value = '➖ Neutre'

$ws.Cells.Item(36, 1).Value = 'SOCIETE GENERALE COTE D''IVOIRE (SGBC)'
$ws.Cells.Item(36, 2).Value = 0
$ws.Cells.Item(36, 3).Value = 1
$ws.Cells.Item(36, 4).Value = -1.85
$ws.Cells.Item(36, 5).Value = -1.85
$ws.Cells.Item(36, 6).Value = '🟡 Observer'
$ws.Cells.Item(36, 7).Value = '➖ Neutre'

$ws.Cells.Item(37, 1).Value = 'TOTALENERGIES MARKETING CI (TTLC)'
$ws.Cells.Item(37, 2).Value = 0
$ws.Cells.Item(37, 3).Value = 1
$ws.Cells.Item(37, 4).Value = -2
$ws.Cells.Item(37, 5).Value = -2
$ws.Cells.Item(37, 6).Value = '🟡 Observer'
$ws.Cells.Item(37, 7).Value = '➖ Neutre'

$ws.Cells.Item(38, 1).Value = 'LOTERIE NATIONALE DU BENIN (LNBB)'
$ws.Cells.Item(38, 2).Value = 0
$ws.Cells.Item(38, 3).Value = 1
$ws.Cells.Item(38, 4).Value = -2.13
$ws.Cells.Item(38, 5).Value = -2.13
$ws.Cells.Item(38, 6).Value = '🟡 Observer'
$ws.Cells.Item(38, 7).Value = '➖ Neutre'

$ws.Cells.Item(39, 1).Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws.Cells.Item(39, 2).Value = 0
$ws.Cells.Item(39, 3).Value = 1
$ws.Cells.Item(39, 4).Value = -2.23
$ws.Cells.Item(39, 5).Value = -2.23
$ws.Cells.Item(39, 6).Value = '🟡 Observer'
$ws.Cells.Item(39, 7).Value = '➖ Neutre'

$ws.Cells.Item(40, 1).Value = 'ORANGE COTE D''IVOIRE (ORAC)'
$ws.Cells.Item(40, 2).Value = 0
$ws.Cells.Item(40, 3).Value = 1
$ws.Cells.Item(40, 4).Value = -3.34
$ws.Cells.Item(40, 5).Value = -3.34
$ws.Cells.Item(40, 6).Value = '🟡 Observer'
$ws.Cells.Item(40, 7).Value = '➖ Neutre'

$ws.Cells.Item(41, 1).Value = 'NEI-CEDA CI (NEIC)'
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(41, 3).Value = 1
$ws.Cells.Item(41, 4).Value = -3.65
$ws.Cells.Item(41, 5).Value = -3.65
$ws.Cells.Item(41, 6).Value = '🟡 Observer'
$ws.Cells.Item(41, 7).Value = '➖ Neutre'

$ws.Cells.Item(42, 1).Value = 'ORAGROUP TOGO (ORGT)'
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(42, 3).Value = 1
$ws.Cells.Item(42, 4).Value = -7.38
$ws.Cells.Item(42, 5).Value = -7.38
$ws.Cells.Item(42, 6).Value = '🟡 Observer'
$ws.Cells.Item(42, 7).Value = '➖ Neutre'

# --- Remove obsolete trailing rows 43-45 (table shrank from 45 to 42 rows) ---
$ws.Rows.Item(43).Delete()
$ws.Rows.Item(43).Delete()
$ws.Rows.Item(43).Delete()

# --- Top_YTD sheet: refresh Progression YTD (%) values ---
$ws2 = $wb.Worksheets.Item("Top_YTD")
$ws2.Cells.Item(2, 2).Value = 506638.76
$ws2.Cells.Item(3, 2).Value = 46423.1
$ws2.Cells.Item(4, 2).Value = 41768.75
$ws2.Cells.Item(5, 2).Value = 34564.89
$ws2.Cells.Item(6, 2).Value = 18635.52
$ws2.Cells.Item(7, 2).Value = 10666.3
$ws2.Cells.Item(8, 2).Value = 9698.870000000001
$ws2.Cells.Item(9, 2).Value = 1625.28
$ws2.Cells.Item(10, 2).Value = 1231.76
$ws2.Cells.Item(11, 2).Value = 1205.7
